$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at L (12). This shifts SIDEBAR_SUBMENU/KODE_JENIS_MANFAAT/
# NAMA_JENIS_MANFAAT headers (L1:N1 -> M1:O1) and the "Setup Jenis
# Manfaat"/"DSBA" data (L2:M2 -> M2:N2) one column to the right. The newly
# inserted L2 cell automatically inherits its style from the K2 cell to its
# left (quote-prefixed style), matching the target style="5".
$ws.Columns(12).Insert()

# Row 1 (headers): the sheet keeps "SIDEBAR_SUBMENU" in column L and puts
# the brand-new "SIDEBAR_SUBMENU_SUBMENU" header in column M, i.e. the
# header that got shifted into M1 by the insert above needs to move back
# to L1, and the inserted (now blank) M1 gets the new header text.
$ws.Cells.Item(1, 13).Copy($ws.Cells.Item(1, 12))
$ws.Cells.Item(1, 13).Value = "SIDEBAR_SUBMENU_SUBMENU"

# Row 2 (data): the newly inserted, blank L2 cell gets the new submenu
# text; the "Setup Jenis Manfaat" value that the insert shifted into M2
# stays there as the new sub-submenu value. Re-apply K2's format (the
# insert already copied it, but setting .Value resets the cell style) so
# L2 keeps the quote-prefixed style used by its left neighbour (K2).
$ws.Cells.Item(2, 12).Value = "Setup Kelengkapan Kepesertaan"
$ws.Cells.Item(2, 11).Copy()
$ws.Cells.Item(2, 12).PasteSpecial(-4122)

# Give the new column a slightly narrower, non-bestFit width.
$ws.Columns(12).ColumnWidth = 15

# Restore the cursor position as recorded after the edit.
$ws.Range("J12").Select()
